$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 16 de Octubre de 2020 a las 10:41"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 8216906
$ws.Range("C4").Value = 591
$ws.Range("D4").Value = 5320386
$ws.Range("E4").Value = 2673793
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 222727

# Row 5: India
$ws.Range("A5").Value = "India"
$ws.Range("B5").Value = 7372394
$ws.Range("C5").Value = 6885
$ws.Range("D5").Value = 6453054
$ws.Range("E5").Value = 807126
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 68
$ws.Range("H5").Value = 112214

# Row 7: Rusia
$ws.Range("A7").Value = "Rusia"
$ws.Range("B7").Value = 1369313
$ws.Range("C7").Value = 15150
$ws.Range("D7").Value = 1056582
$ws.Range("E7").Value = 289008
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 232
$ws.Range("H7").Value = 23723

# Row 21: Filipinas
$ws.Range("A21").Value = "Filipinas"
$ws.Range("B21").Value = 351750
$ws.Range("C21").Value = 3139
$ws.Range("D21").Value = 294865
$ws.Range("E21").Value = 50354
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 34
$ws.Range("H21").Value = 6531

# Row 22: Indonesia
$ws.Range("A22").Value = "Indonesia"
$ws.Range("B22").Value = 349160
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 273661
$ws.Range("E22").Value = 63231
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 12268

# Row 23: Alemania
$ws.Range("A23").Value = "Alemania"
$ws.Range("B23").Value = 348816
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 284600
$ws.Range("E23").Value = 54406
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 9810

# Row 35: Polonia
$ws.Range("A35").Value = "Polonia"
$ws.Range("B35").Value = 149903
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 87773
$ws.Range("E35").Value = 58822
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 3308

# Row 65: Singapur
$ws.Range("A65").Value = "Singapur"
$ws.Range("B65").Value = 57901
$ws.Range("C65").Value = 9
$ws.Range("D65").Value = 57764
$ws.Range("E65").Value = 109
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 28

# Row 86: Eslovaquia
$ws.Range("A86").Value = "Eslovaquia"
$ws.Range("B86").Value = 26300
$ws.Range("C86").Value = 2075
$ws.Range("D86").Value = 7182
$ws.Range("E86").Value = 19047
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 71

# Row 87: Corea del Sur
$ws.Range("A87").Value = "Corea del Sur"
$ws.Range("B87").Value = 25035
$ws.Range("C87").Value = 47
$ws.Range("D87").Value = 23180
$ws.Range("E87").Value = 1414
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 2
$ws.Range("H87").Value = 441

# Row 89: Croacia
$ws.Range("A89").Value = "Croacia"
$ws.Range("B89").Value = 23665
$ws.Range("C89").Value = 1131
$ws.Range("D89").Value = 19087
$ws.Range("E89").Value = 4233
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 1
$ws.Range("H89").Value = 345

# Row 120: Lituania
$ws.Range("A120").Value = "Lituania"
$ws.Range("B120").Value = 7041
$ws.Range("C120").Value = 281
$ws.Range("D120").Value = 3035
$ws.Range("E120").Value = 3894
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 2
$ws.Range("H120").Value = 112

# Row 121: Guadalupe
$ws.Range("A121").Value = "Guadalupe"
$ws.Range("B121").Value = 6908
$ws.Range("C121").Value = 0
$ws.Range("D121").Value = 2199
$ws.Range("E121").Value = 4613
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 96

# Row 141: Estonia
$ws.Range("A141").Value = "Estonia"
$ws.Range("B141").Value = 4017
$ws.Range("C141").Value = 37
$ws.Range("D141").Value = 3137
$ws.Range("E141").Value = 812
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 68

# Row 150: Letonia
$ws.Range("A150").Value = "Letonia"
$ws.Range("B150").Value = 3204
$ws.Range("C150").Value = 148
$ws.Range("D150").Value = 1329
$ws.Range("E150").Value = 1833
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 1
$ws.Range("H150").Value = 42

# Row 151: Principado de Andorra
$ws.Range("A151").Value = "Principado de Andorra"
$ws.Range("B151").Value = 3190
$ws.Range("C151").Value = 0
$ws.Range("D151").Value = 2011
$ws.Range("E151").Value = 1120
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 59
